$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("DPLKINV083-001")
$ws2 = $wb.Worksheets.Item("DPLKINV083-002")
$ws3 = $wb.Worksheets.Item("DPLKINV083-003")

# Sheet 1: DPLKINV083-001 - update No. Urut from 1369 to 2962
$ws1.Range("L2").Value = 2962
$ws1.Range("F2").Value = "Username : 32070;`nPassword : bni1234;`nRole : 18/19 - Pimpinan Kelompok Investasi/Pengelola Investasi;`nNo. Urut : 2962;`nStatus Verifikasi : 1 : Setuju;`nKeterangan Verifikasi : DATA APPROVAL"

# Sheet 2: DPLKINV083-002 - update No. Urut from 1369 to 2962
$ws2.Range("L2").Value = 2962
$ws2.Range("F2").Value = "Username : 32070;`nPassword : bni1234;`nRole : 18/19 - Pimpinan Kelompok Investasi/Pengelola Investasi;`nNo. Urut : 2962;`nStatus Verifikasi : 2 : Batalkan Verifikasi;`nKeterangan Verifikasi : DATA AKAN DIHAPUS"

# Sheet 3: DPLKINV083-003 - update No. Urut from 1369 to 2962
$ws3.Range("L2").Value = 2962
$ws3.Range("F2").Value = "Username : 32070;`nPassword : bni1234;`nRole : 18/19 - Pimpinan Kelompok Investasi/Pengelola Investasi;`nNo. Urut : 2962;`nStatus Verifikasi : 0 : Kembalikan ke Register;`nKeterangan Verifikasi : DATA AKAN DIKEMBALIKAN UNTUK DIEDIT"

# Update sheet views (selection / scroll position) to match the target state
$ws1.Activate()
$ws1.Application.ActiveWindow.ScrollColumn = 1
$ws1.Range("G2").Select()

$ws2.Activate()
$ws2.Range("G2").Select()

$ws3.Activate()
$ws3.Application.ActiveWindow.ScrollColumn = 6
$ws3.Range("Q2").Select()
